# -----------------------------------------------------------------------
# Adapt column header formatting to respective input file names.
#   *_old  -> *_FV2310
#   *_new  -> *_FV2404
# Wrap the header row (and full data range) in an Excel Table ("Table1")
# and freeze the header row.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the 21 header cells in row 1 (A1:U1). Writing new literal
#    values here updates the shared-string table under the hood.
$newHeaders = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeaders[$i]
}

# 2. Turn A1:U71 into an Excel Table named "Table1".
#    Stash the header row's existing manual formatting (bold, grey fill,
#    borders, centered+wrapped) in a scratch cell, clear the header format
#    before creating the table (so Excel doesn't bake a header dxf /
#    default table style into the table definition), then restore the
#    formatting from the stash afterwards.
$stash = $ws.Range("W1")
$headerRange = $ws.Range("A1:U1")

$ws.Range("A1").Copy() | Out-Null
$stash.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$headerRange.ClearFormats()

$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U71"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

$stash.Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$stash.Clear() | Out-Null

# 3. Freeze the header row (split after row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
